# Weekly data refresh: insert the new week's observation for Zapallo
# italiano at Vega Monumental Concepcion ahead of the previous one, pushing
# the prior row (and everything below) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift current row 113 (and below) down to make room for the new entry.
# This also duplicates the formatting (incl. the date-format style on
# column D) onto the new row 114.
$ws.Rows.Item(113).Insert()

# Populate the new week's values into row 113.
$ws.Range("A113").Value = 11
$ws.Range("B113").Value = "Vega Monumental Concepción"
$ws.Range("C113").Value = "Bíobío"
$ws.Range("D113").Value = 44656
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = 100112032
$ws.Range("G113").Value = "Zapallo italiano"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 220
$ws.Range("K113").Value = 8000
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = 8455
$ws.Range("N113").Value = "$/caja 50 unidades"
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 169
$ws.Range("Q113").Value = 50
$ws.Range("R113").Value = "Hortaliza"
